$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.470.16'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '3.333.56'
$ws.Range("E3").Value = '  +7.24%  '
$ws.Range("E4").Value = '  +1.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.64'
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.61'
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("E7").Value = '  +1.34%  '
$ws.Range("D8").Value = '3.181.47'
$ws.Range("E8").Value = '  +2.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("E11").Value = '  +6.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.466'
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.03'
$ws.Range("E14").Value = '  +2.66%  '
$ws.Range("D15").Value = '3.880.50'
$ws.Range("E15").Value = '  +7.15%  '
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("E17").Value = '  +4.13%  '
$ws.Range("D18").Value = '3.221.42'
$ws.Range("E18").Value = '  +3.62%  '
$ws.Range("D19").Value = '64.419.24'
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.56'
$ws.Range("E20").Value = '  +2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.73'
$ws.Range("E21").Value = '  +3.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.754'
$ws.Range("E22").Value = '  +3.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.77'
$ws.Range("E23").Value = '  +4.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").Value = '  +11.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.66'
$ws.Range("E25").Value = '  +5.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.49'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("E27").Value = '  +8.55%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  +2.43%  '
$ws.Range("E30").Value = '  +2.73%  '
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.32'
$ws.Range("E32").Value = '  +4.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.122'
$ws.Range("E33").Value = '  +10.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.51'
$ws.Range("E34").Value = '  +6.47%  '
$ws.Range("D35").Value = '0.0₃0861'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("E36").Value = '  +3.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.29'
$ws.Range("E37").Value = '  +4.67%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.34'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '471.11'
$ws.Range("E40").Value = '  +7.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.45'
$ws.Range("E41").Value = '  +8.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '51.77'
$ws.Range("E42").Value = '  +2.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.301'
$ws.Range("E43").Value = '  +8.70%  '
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("D45").Value = '2.946.41'
$ws.Range("E45").Value = '  +2.18%  '
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.41'
$ws.Range("E47").Value = '  +7.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.80'
$ws.Range("E48").Value = '  +5.59%  '
$ws.Range("E49").Value = '  +6.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.07'
$ws.Range("E51").Value = '  +3.54%  '
